$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix a typo in an existing entry: 12-12-2017 -> 12-28-2017
$ws.Range("A13").Value = "1:54PM 12-28-2017"

# The "Total Project Hours:" summary row (previously row 17) is being
# replaced by new data rows, and the summary moves down to row 21
# (row 20 is left blank).
$ws.Range("A17").Value = "3:04 PM 12-31-2017"
$ws.Range("B17").Value = "4:47 PM 12-31-2017"
$ws.Range("C17").Value = 103

$ws.Range("A18").Value = "6:23 PM 12-31-2017"
$ws.Range("B18").Value = "8:19 PM 12-31-2017"
$ws.Range("C18").Value = 116

$ws.Range("A19").Value = "8:32PM 12-31-2017"
$ws.Range("B19").Value = "10:32PM 12-31-2017"
$ws.Range("C19").Value = 120

$ws.Range("A21").Value = "Total Project Hours:"
$ws.Range("C21").Formula = "=SUM(C2:C20)/60"

# Match the number format used for comparable time-of-day strings
# elsewhere in column A/B (e.g. B8) for the new summary anchor cell A17.
$ws.Range("A17").NumberFormat = "HH:MM:SS\ AM/PM"

$ws.Range("B19").Select() | Out-Null
